# Apply "changes up through module 3 based on ruth feedback"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell content updates (Assignments column, C) ---

# C21: Intro to Scientific Writing -> add Writing Fellow Appointment
$ws.Range("C21").Value = "Intro to Scientific Writing/Writing Fellow Appointment"

# C33: Project 2 Peer Review -> add Writing Fellow Appointment
$ws.Range("C33").Value = "Project 2 Peer Review/Writing Fellow Appointment"

# C39: "Final Project Update" moved down into C40 -> clear this cell
$ws.Range("C39").Value = ""

# C40: Quiz 5 -> merge with Final Project Update
$ws.Range("C40").Value = "Quiz 5/Final Project Update"

# C41: was empty -> new item
$ws.Range("C41").Value = "Regression in the real world"

# C44: Final Project Papers - Penultimate Draft -> two dated lines, now wrapped
$ws.Range("C44").Value = "(23-Apr) Final Project Papers - Penultimate Draft`n(24-Apr) Group Member Evaluation "
$ws.Range("C44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 30.5

# B51: footnote time changed from 1:15pm to 1:00pm
$ws.Range("B51").Value = "* - Assignment due before class (1:00pm)"

# --- Column width ---
$ws.Columns.Item(3).ColumnWidth = 47.45

# --- View state (zoom / selection) ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("F28").Select()
